# Update the "Förändrad" (changed) date column (C) for every data row,
# and add a friendly display-text second argument to the HYPERLINK()
# formulas in columns S, T, V, W, X, Y for the rows that have them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastDataRow = 171
$lastLinkRow = 13

$linkCols = @(
    @{ Col = "S"; Folder = "artfynd";        Ext = ".xlsx" },
    @{ Col = "T"; Folder = "kartor";         Ext = ".png"  },
    @{ Col = "V"; Folder = "klagomål";       Ext = ".docx" },
    @{ Col = "W"; Folder = "klagomålsmail";  Ext = ".docx" },
    @{ Col = "X"; Folder = "tillsyn";        Ext = ".docx" },
    @{ Col = "Y"; Folder = "tillsynsmail";   Ext = ".docx" }
)

for ($row = 2; $row -le $lastDataRow; $row++) {

    # Column C: "Förändrad" date serial moves from 45184 to 45186.
    $ws.Cells.Item($row, 3).Value = 45186

    if ($row -le $lastLinkRow) {
        $id = $ws.Range("A$row").Text

        foreach ($link in $linkCols) {
            $url = "https://klasma.github.io/Logging_GOTENE/" + $link.Folder + "/" + $id + $link.Ext
            $formula = '=HYPERLINK("' + $url + '", "' + $id + '")'
            $ws.Range($link.Col + "$row").Formula = $formula
        }
    }
}
